# COVID-19 Bangladesh DataSheet update (12 June 2020 data point)
# Appends a new row (row 98, date 2020-06-12 / serial 43994) to each of the
# three sheets (Confirmed, Recoverd, Death), carrying the "running total"
# formula down one more row and leaving the new row selected, mirroring
# what a user does in Excel: select the last data row, copy it down, then
# type in the new day's date and new-case count.

$wb = $excel.ActiveWorkbook

function Add-DailyRow {
    param(
        [string]$SheetName,
        [int]$DateSerial,
        [int]$NewCount
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Bring formatting (date format in col A, centered number format in
    # col B/C) down onto the new row by copying the last existing row's
    # formats, the same way Excel does when you drag the fill handle down.
    $ws.Range("A97:C97").Copy()
    $ws.Range("A98:C98").PasteSpecial(-4122)

    # New date value.
    $ws.Range("A98").Value = $DateSerial

    # Running-total formula, same pattern as every other row: previous
    # day's total plus today's new count.
    $ws.Range("B98").Formula = "=SUM(B97+C98)"

    # Today's new-case (or new-recovered / new-death) count.
    $ws.Range("C98").Value = $NewCount

    # Match Excel's habit of leaving the two most-recently-touched cells
    # selected after the edit.
    $ws.Activate()
    $ws.Range("B97:B98").Select()
}

Add-DailyRow "Confirmed" 43994 3471
Add-DailyRow "Recoverd"  43994 502
Add-DailyRow "Death"     43994 46
